$d = $word.ActiveDocument

# The paragraph "{m:Sequence{''.sampleTable(),''.sampleTable(),''.sampleTable()}}"
# currently has its leading "{m" in a single run and its trailing "}}" in a
# single run. The TokenIteratorFieldRewriterSplit parser expects the field
# delimiters ("{" for the opening brace and the final "}" for the closing
# brace) to live in their own runs, so we split those two runs in place
# without touching the text or any other run's formatting.
#
# A collapsed Bookmark.Add/Delete pair at the desired character offset is
# used purely as a mechanism to force Word to break the run at that exact
# position -- it leaves no trace (no bookmark, no extra run formatting) in
# the saved document.

function Split-RunAt($position) {
    $bookmarkName = "__runsplit__"
    $d.Bookmarks.Add($bookmarkName, $d.Range($position, $position)) | Out-Null
    $d.Bookmarks($bookmarkName).Delete()
}

# --- Split "{m" into "{" and "m" ---
$openBrace = $d.Content
$openBrace.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-RunAt ($openBrace.Start + 1)

# --- Split the trailing "}}" into "}" and "}" ---
$closeBraces = $d.Content
$closeBraces.Find.Execute("}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Split-RunAt ($closeBraces.End - 1)
